$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header labels: "_old" suffix -> "_FV2210", "_new" suffix -> "_FV2304" ---
# Columns A-J (1-10) held the "*_old" headers, column K (11) is "diff",
# columns L-U (12-21) held the "*_new" headers.
$fv2210Headers = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210"
)
$fv2304Headers = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt $fv2210Headers.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $fv2210Headers[$i]
}
for ($i = 0; $i -lt $fv2304Headers.Count; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = $fv2304Headers[$i]
}

# --- Turn the data range into an Excel table (ListObject) ---
$rng = $ws.Range("A1:U73")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"

# --- Freeze the header row (pane split after row 1) ---
$ws.Activate()
[void]$ws.Range("A2").Select()
$appWin = $excel.ActiveWindow
$appWin.FreezePanes = $true
